$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate formatting of the last existing data row (row 7) down through the new rows (8-15)
# so the newly added rows visually match the existing bordered data rows.
$ws.Range("A7:H7").Copy() | Out-Null
$ws.Range("A8:H15").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = 0

# Update existing rows (2-7) with the refreshed client/candidate data, and
# populate the newly appended rows (8-15) with the additional candidates.
$ws.Range("A2").Value = "yPFIq555"
$ws.Range("B2").Value = 231011227
$ws.Range("C2").Value = "bsnwroa36"
$ws.Range("D2").Value = "A6%Fuj`$2"
$ws.Range("E2").Value = "MR"
$ws.Range("F2").Value = "gYVZaiqa"
$ws.Range("G2").Value = "bvVM"
$ws.Range("H2").Value = "Candidate"

$ws.Range("A3").Value = "oPMcI828"
$ws.Range("B3").Value = 231011226
$ws.Range("C3").Value = "naljpwc93"
$ws.Range("D3").Value = "k35hM#!H"
$ws.Range("E3").Value = "MR"
$ws.Range("F3").Value = "zGfIUZZs"
$ws.Range("G3").Value = "zkzt"
$ws.Range("H3").Value = "Candidate"

$ws.Range("A4").Value = "AQYIj613"
$ws.Range("B4").Value = 231011225
$ws.Range("C4").Value = "hzefdef87"
$ws.Range("D4").Value = "u%6y5AR&"
$ws.Range("E4").Value = "MR"
$ws.Range("F4").Value = "UOOwvTOZ"
$ws.Range("G4").Value = "LYdT"
$ws.Range("H4").Value = "Candidate"

$ws.Range("A5").Value = "RJYmE208"
$ws.Range("B5").Value = 231011224
$ws.Range("C5").Value = "hlxlque29"
$ws.Range("D5").Value = "ZR&28!mc"
$ws.Range("E5").Value = "MR"
$ws.Range("F5").Value = "QWYFBbIw"
$ws.Range("G5").Value = "HkFJ"
$ws.Range("H5").Value = "Candidate"

$ws.Range("A6").Value = "DjdcV735"
$ws.Range("B6").Value = 231011223
$ws.Range("C6").Value = "djgjijf84"
$ws.Range("D6").Value = "r!34Z&My"
$ws.Range("E6").Value = "MR"
$ws.Range("F6").Value = "BmaxKoiX"
$ws.Range("G6").Value = "OOwm"
$ws.Range("H6").Value = "Candidate"

$ws.Range("A7").Value = "GNkKw909"
$ws.Range("B7").Value = 231011222
$ws.Range("C7").Value = "avxdomr45"
$ws.Range("D7").Value = "Wa45`$dA!"
$ws.Range("E7").Value = "MR"
$ws.Range("F7").Value = "vYDSrtYs"
$ws.Range("G7").Value = "jLbh"
$ws.Range("H7").Value = "Candidate"

$ws.Range("A8").Value = "MNCDF699"
$ws.Range("B8").Value = 231011221
$ws.Range("C8").Value = "wkwmqvp50"
$ws.Range("D8").Value = "k4!b#5FA"
$ws.Range("E8").Value = "MR"
$ws.Range("F8").Value = "DwYDUTJm"
$ws.Range("G8").Value = "czJN"
$ws.Range("H8").Value = "Candidate"

$ws.Range("A9").Value = "HEDvZ511"
$ws.Range("B9").Value = 231011220
$ws.Range("C9").Value = "xfpirxx37"
$ws.Range("D9").Value = "n&6u#D2C"
$ws.Range("E9").Value = "MR"
$ws.Range("F9").Value = "gKPSQbFj"
$ws.Range("G9").Value = "qDuu"
$ws.Range("H9").Value = "Candidate"

$ws.Range("A10").Value = "udcoU491"
$ws.Range("B10").Value = 231011219
$ws.Range("C10").Value = "zmhqktu15"
$ws.Range("D10").Value = "J`$Y!e47q"
$ws.Range("E10").Value = "MR"
$ws.Range("F10").Value = "wHVKUpqk"
$ws.Range("G10").Value = "RyKU"
$ws.Range("H10").Value = "Candidate"

$ws.Range("A11").Value = "YqcqJ884"
$ws.Range("B11").Value = 231011218
$ws.Range("C11").Value = "bjtuchu21"
$ws.Range("D11").Value = "d&2!SCz8"
$ws.Range("E11").Value = "MR"
$ws.Range("F11").Value = "MWKHhhlj"
$ws.Range("G11").Value = "Eiuz"
$ws.Range("H11").Value = "Candidate"

$ws.Range("A12").Value = "RgtQS935"
$ws.Range("B12").Value = 231011217
$ws.Range("C12").Value = "owgsftg60"
$ws.Range("D12").Value = "zm5X4#&N"
$ws.Range("E12").Value = "MR"
$ws.Range("F12").Value = "eAvsguvE"
$ws.Range("G12").Value = "fbpn"
$ws.Range("H12").Value = "Candidate"

$ws.Range("A13").Value = "smwan615"
$ws.Range("B13").Value = 231011216
$ws.Range("C13").Value = "znxshqo57"
$ws.Range("D13").Value = "rGy9Q#`$4"
$ws.Range("E13").Value = "MR"
$ws.Range("F13").Value = "IPuypWoz"
$ws.Range("G13").Value = "fIDm"
$ws.Range("H13").Value = "Candidate"

$ws.Range("A14").Value = "ptSuC615"
$ws.Range("B14").Value = 231011215
$ws.Range("C14").Value = "isaeobj36"
$ws.Range("D14").Value = "U95#pZ`$a"
$ws.Range("E14").Value = "MR"
$ws.Range("F14").Value = "wyXnLvVP"
$ws.Range("G14").Value = "MwLH"
$ws.Range("H14").Value = "Candidate"

$ws.Range("A15").Value = "QLxLZ247"
$ws.Range("B15").Value = 231011214
$ws.Range("C15").Value = "mnqwsxj76"
$ws.Range("D15").Value = "wN65Xt&`$"
$ws.Range("E15").Value = "MR"
$ws.Range("F15").Value = "tfMURMTW"
$ws.Range("G15").Value = "yEYa"
$ws.Range("H15").Value = "Candidate"

# Refresh the selection to cover the full populated range
$ws.Range("A1:H15").Select() | Out-Null